$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038885193359007
$ws.Range("D2").Value = 1.047474620495898
$ws.Range("E2").Value = 1.037486875509379
$ws.Range("F2").Value = 1.054512006590757
$ws.Range("I2").Value = 1.042805451800082
$ws.Range("J2").Value = 1.04398020588261
$ws.Range("K2").Value = 1.050237156097804
$ws.Range("L2").Value = 1.040277620747508
$ws.Range("M2").Value = 1.057255013447404

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039795956732021
$ws.Range("D3").Value = 1.048222821768606
$ws.Range("E3").Value = 1.03826045660304
$ws.Range("F3").Value = 1.05543656121392
$ws.Range("I3").Value = 1.04307085780012
$ws.Range("J3").Value = 1.044536140879971
$ws.Range("K3").Value = 1.050797373964917
$ws.Range("L3").Value = 1.040861115461958
$ws.Range("M3").Value = 1.057992541464523

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.040385722909166
$ws.Range("D4").Value = 1.04870729900636
$ws.Range("E4").Value = 1.038761763015982
$ws.Range("F4").Value = 1.056035616432426
$ws.Range("I4").Value = 1.043241465887922
$ws.Range("J4").Value = 1.04489567114256
$ws.Range("K4").Value = 1.051159547367155
$ws.Range("L4").Value = 1.041238763503105
$ws.Range("M4").Value = 1.058469947131738

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040633765116262
$ws.Range("D5").Value = 1.04891105341903
$ws.Range("E5").Value = 1.038972689765685
$ws.Range("F5").Value = 1.05628765084645
$ws.Range("I5").Value = 1.043312918977299
$ws.Range("J5").Value = 1.045046769503286
$ws.Range("K5").Value = 1.051311726152447
$ws.Range("L5").Value = 1.041397546672969
$ws.Range("M5").Value = 1.058670689250998

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040675418588942
$ws.Range("D6").Value = 1.04894526934122
$ws.Range("E6").Value = 1.03900811566296
$ws.Range("F6").Value = 1.056329979725796
$ws.Range("I6").Value = 1.043324900380952
$ws.Range("J6").Value = 1.045072136728648
$ws.Range("K6").Value = 1.051337272980572
$ws.Range("L6").Value = 1.041424208202201
$ws.Range("M6").Value = 1.058704397089772

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040389036852083
$ws.Range("D7").Value = 1.048710021269555
$ws.Range("E7").Value = 1.038764580733797
$ws.Range("F7").Value = 1.056038983377611
$ws.Range("I7").Value = 1.043242421711505
$ws.Range("J7").Value = 1.04489769031679
$ws.Range("K7").Value = 1.051161581098988
$ws.Range("L7").Value = 1.041240885094674
$ws.Range("M7").Value = 1.058472629299155

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.039192898017556
$ws.Range("D8").Value = 1.047727407144432
$ws.Range("E8").Value = 1.037748155157599
$ws.Range("F8").Value = 1.054824296255552
$ws.Range("I8").Value = 1.042895380024912
$ws.Range("J8").Value = 1.044168126856006
$ws.Range("K8").Value = 1.050426550982895
$ws.Range("L8").Value = 1.040474796635115
$ws.Range("M8").Value = 1.057504227332396

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03708858019922
$ws.Range("D9").Value = 1.045998591054991
$ws.Range("E9").Value = 1.035962875018762
$ws.Range("F9").Value = 1.052690105431475
$ws.Range("I9").Value = 1.042275243440337
$ws.Range("J9").Value = 1.042881080812116
$ws.Range("K9").Value = 1.049128893447743
$ws.Range("L9").Value = 1.039125573907679
$ws.Range("M9").Value = 1.055799182211702

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035688076950097
$ws.Range("D10").Value = 1.044847931657325
$ws.Range("E10").Value = 1.034776669358332
$ws.Range("F10").Value = 1.051271597099869
$ws.Range("I10").Value = 1.041856072786135
$ws.Range("J10").Value = 1.042022128161583
$ws.Range("K10").Value = 1.048262212642426
$ws.Range("L10").Value = 1.038226644450435
$ws.Range("M10").Value = 1.054663504951023

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035082220986127
$ws.Range("D11").Value = 1.044350148371967
$ws.Range("E11").Value = 1.034263991724344
$ws.Range("F11").Value = 1.050658402151112
$ws.Range("I11").Value = 1.041673212433681
$ws.Range("J11").Value = 1.041649986192983
$ws.Range("K11").Value = 1.047886570081208
$ws.Range("L11").Value = 1.037837542904706
$ws.Range("M11").Value = 1.054172001844733

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034857266005836
$ws.Range("D12").Value = 1.044165320004256
$ws.Range("E12").Value = 1.034073705699175
$ws.Range("F12").Value = 1.050430790222488
$ws.Range("I12").Value = 1.041605086749761
$ws.Range("J12").Value = 1.041511725570838
$ws.Range("K12").Value = 1.047746986076883
$ws.Range("L12").Value = 1.037693035519498
$ws.Range("M12").Value = 1.053989474915377

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034905515662821
$ws.Range("D13").Value = 1.044204963103031
$ws.Range("E13").Value = 1.034114516100081
$ws.Range("F13").Value = 1.050479606672133
$ws.Range("I13").Value = 1.041619709125786
$ws.Range("J13").Value = 1.041541384308975
$ws.Range("K13").Value = 1.047776929729502
$ws.Range("L13").Value = 1.037724031819579
$ws.Range("M13").Value = 1.054028625766946

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035063624358239
$ws.Range("D14").Value = 1.044334868954852
$ws.Range("E14").Value = 1.034248259650136
$ws.Range("F14").Value = 1.050639584487771
$ws.Range("I14").Value = 1.041667585285164
$ws.Range("J14").Value = 1.041638558143776
$ws.Range("K14").Value = 1.047875033111783
$ws.Range("L14").Value = 1.037825597418871
$ws.Range("M14").Value = 1.054156913304727

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035161051891461
$ws.Range("D15").Value = 1.044414917623224
$ws.Range("E15").Value = 1.034330682759431
$ws.Range("F15").Value = 1.05073817281203
$ws.Range("I15").Value = 1.041697056458752
$ws.Range("J15").Value = 1.041698426135414
$ws.Range("K15").Value = 1.047935470768542
$ws.Range("L15").Value = 1.037888178318228
$ws.Range("M15").Value = 1.054235960723754

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035728297741929
$ws.Range("D16").Value = 1.04488097770419
$ws.Range("E16").Value = 1.034810714377112
$ws.Range("F16").Value = 1.051312314622947
$ws.Range("I16").Value = 1.041868180094455
$ws.Range("J16").Value = 1.042046821663031
$ws.Range("K16").Value = 1.048287135222619
$ws.Range("L16").Value = 1.038252470899272
$ws.Range("M16").Value = 1.054696129823166

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036084269765519
$ws.Range("D17").Value = 1.045173448994315
$ws.Range("E17").Value = 1.035112082841035
$ws.Range("F17").Value = 1.051672735002499
$ws.Range("I17").Value = 1.041975158642328
$ws.Range("J17").Value = 1.042265305505037
$ws.Range("K17").Value = 1.048507628326747
$ws.Range("L17").Value = 1.038481020490178
$ws.Range("M17").Value = 1.054984850264516

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.036291957262284
$ws.Range("D18").Value = 1.045344086820764
$ws.Range("E18").Value = 1.035287958229672
$ws.Range("F18").Value = 1.051883061251677
$ws.Range("I18").Value = 1.042037426476424
$ws.Range("J18").Value = 1.042392723123917
$ws.Range("K18").Value = 1.048636203026717
$ws.Range("L18").Value = 1.038614343115575
$ws.Range("M18").Value = 1.055153280305451

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.036362782633957
$ws.Range("D19").Value = 1.045402277375716
$ws.Range("E19").Value = 1.035347942797044
$ws.Range("F19").Value = 1.0519547938389
$ws.Range("I19").Value = 1.042058635969661
$ws.Range("J19").Value = 1.042436165764052
$ws.Range("K19").Value = 1.048680037674345
$ws.Range("L19").Value = 1.038659804944707
$ws.Range("M19").Value = 1.055210714679108

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036046071635564
$ws.Range("D20").Value = 1.045142065003332
$ws.Range("E20").Value = 1.035079739303295
$ws.Range("F20").Value = 1.05163405504176
$ws.Range("I20").Value = 1.041963694396128
$ws.Range("J20").Value = 1.042241866350027
$ws.Range("K20").Value = 1.048483975136319
$ws.Range("L20").Value = 1.038456497876133
$ws.Range("M20").Value = 1.054953870771862

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035017062880568
$ws.Range("D21").Value = 1.044296612961496
$ws.Range("E21").Value = 1.034208871486091
$ws.Range("F21").Value = 1.050592470690967
$ws.Range("I21").Value = 1.041653492554026
$ws.Range("J21").Value = 1.041609943706073
$ws.Range("K21").Value = 1.047846145582814
$ws.Range("L21").Value = 1.03779568825824
$ws.Range("M21").Value = 1.054119134729925

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034370587555381
$ws.Range("D22").Value = 1.043765452350917
$ws.Range("E22").Value = 1.033662163563671
$ws.Range("F22").Value = 1.049938488430096
$ws.Range("I22").Value = 1.041457280958862
$ws.Range("J22").Value = 1.041212452741867
$ws.Range("K22").Value = 1.047444807058329
$ws.Range("L22").Value = 1.037380340318344
$ws.Range("M22").Value = 1.053594530518802

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034713248284116
$ws.Range("D23").Value = 1.04404699143019
$ws.Range("E23").Value = 1.033951903550387
$ws.Range("F23").Value = 1.050285090705078
$ws.Range("I23").Value = 1.041561407643034
$ws.Range("J23").Value = 1.041423186584053
$ws.Range("K23").Value = 1.04765759321093
$ws.Range("L23").Value = 1.037600511498988
$ws.Range("M23").Value = 1.053872611157134

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036063331556299
$ws.Range("D24").Value = 1.045156245941312
$ws.Range("E24").Value = 1.035094353670877
$ws.Range("F24").Value = 1.051651532544882
$ws.Range("I24").Value = 1.041968875000159
$ws.Range("J24").Value = 1.042252457558098
$ws.Range("K24").Value = 1.048494663103791
$ws.Range("L24").Value = 1.038467578548208
$ws.Range("M24").Value = 1.054967868996161

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037632183130363
$ws.Range("D25").Value = 1.046445205547308
$ws.Range("E25").Value = 1.036423717410374
$ws.Range("F25").Value = 1.053241096009777
$ws.Range("I25").Value = 1.042436579300922
$ws.Range("J25").Value = 1.043213980036389
$ws.Range("K25").Value = 1.04946465153743
$ws.Range("L25").Value = 1.039474287816732
$ws.Range("M25").Value = 1.056239803142023

